# Applies the LOB1003 syllabus-sheet edit:
#  - removes the standalone "5840692 - Diovana Aparecida dos Santos Napoleão"
#    row (old row 13, a label-less row sitting right under "Docentes
#    responsáveis:") which shifts everything below it up by one row
#  - re-purposes several now-orphaned text cells with new content
#  - gives the "Método:" / "Critério:" / "Norma de recuperação:" rows a
#    60pt custom row height (matching the new text lengths)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the old row 13 (only B13/C13 held "5840692 - Diovana Aparecida
#    dos Santos Napoleão"; A13 was empty). Rows 14-22 shift up to 13-21.
$ws.Rows("13:13").Delete()

# 2) Row 10 "Objetivos:" now shows the docente's name instead of the
#    Portuguese "Fornecer fundamentos..." objectives text.
$ws.Range("B10").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"
$ws.Range("C10").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"

# 3) Row 13 "Programa resumido:" (previously row 14) now just says
#    "Semestral" instead of the long summary text.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# 4) Row 15 "Programa:" (previously row 16) now holds a date instead of
#    the full bulleted syllabus text.
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

# 5) Row 18 "Método:" (previously row 19) now shows the docente's name,
#    and gains a 60pt custom row height.
$ws.Range("B18").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"
$ws.Range("C18").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"
$ws.Rows("18:18").RowHeight = 60

# 6) Row 19 "Critério:" (previously row 20) now holds the evaluation
#    method text, with a 60pt custom row height.
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Rows("19:19").RowHeight = 60

# 7) Row 20 "Norma de recuperação:" (previously row 21) now holds the
#    criterion text "NF>= 5,0.", with a 60pt custom row height.
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."
$ws.Rows("20:20").RowHeight = 60

# 8) Row 21 "Bibliografia:" (previously row 22) now holds the recovery
#    norm text instead of the full bibliography list.
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
